$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: classification label for each feedback example.
# Row 1 is the header ("Label"); rows 2-54 are "Positive"/"Negative".
$labels = @("Label","Negative","Positive","Positive","Positive","Negative","Positive","Positive","Negative","Positive","Positive","Negative","Positive","Positive","Negative","Positive","Positive","Positive","Negative","Positive","Positive","Negative","Positive","Positive","Positive","Positive","Negative","Negative","Positive","Positive","Positive","Positive","Negative","Positive","Negative","Positive","Positive","Positive","Positive","Positive","Positive","Positive","Positive","Positive","Positive","Negative","Positive","Positive","Positive","Negative","Positive","Negative","Negative","Positive")

for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $labels[$i]
}

# Turn on AutoFilter across the now-two-column table.
$dataRange = $ws.Range("A1:B54")
$dataRange.AutoFilter() | Out-Null

# Excel records the filter range as a hidden workbook-scoped defined name.
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=" + $ws.Name + "!`$A`$1:`$B`$54")
$filterName.Visible = $false

Write-Output "done"
